$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 53
$ws1.Range("F7").Value = 584
$ws1.Range("F8").Value = 79
$ws1.Range("F9").Value = 8506
$ws1.Range("F10").Value = 795
$ws1.Range("F11").Value = 314
$ws1.Range("F12").Value = 1133
$ws1.Range("F13").Value = 916
$ws1.Range("F14").Value = 83
$ws1.Range("F15").Value = 43
$ws1.Range("F16").Value = 223
$ws1.Range("F17").Value = 194
$ws1.Range("F20").Value = 963

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 53
$ws4.Range("F9").Value = 584
$ws4.Range("F10").Value = 79
$ws4.Range("F11").Value = 8506
$ws4.Range("F12").Value = 795
$ws4.Range("F13").Value = 314
$ws4.Range("F14").Value = 1133
$ws4.Range("F15").Value = 916
$ws4.Range("F16").Value = 83
$ws4.Range("F17").Value = 43
$ws4.Range("F18").Value = 223
$ws4.Range("F19").Value = 194
$ws4.Range("F22").Value = 963
